$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = -5.815899999999998
$ws.Range("A9").Value = -20.33769999999997
$ws.Range("A18").Value = -23.09220000000001
$ws.Range("A20").Value = -22.26920000000003
$ws.Range("E21").Value = 13.0674
